$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 310.4375
$ws.Range("I6").Value = 191.88889
$ws.Range("K6").Value = 575.6666700000001
$ws.Range("M6").Value = -463.6666700000001
$ws.Range("H11").Value = 74.73333
$ws.Range("I11").Value = 74.73333
$ws.Range("K11").Value = 74.73333
$ws.Range("M11").Value = 65.26667
$ws.Range("H12").Value = 299.72726
$ws.Range("I12").Value = 310.8889
$ws.Range("K12").Value = 310.8889
$ws.Range("M12").Value = -140.8889
$ws.Range("H40").Value = 3989
$ws.Range("I40").Value = 3483.5
$ws.Range("J40").Value = 5000
$ws.Range("K40").Value = 3483.5
$ws.Range("L40").Value = 5000
$ws.Range("M40").Value = -3308.5
$ws.Range("N40").Value = -5350
$ws.Range("H41").Value = 25475.25
$ws.Range("J41").Value = 33483.668
$ws.Range("L41").Value = 33483.668
$ws.Range("N41").Value = -34363.668
$ws.Range("H54").Value = 23571.428
$ws.Range("I54").Value = 20000
$ws.Range("K54").Value = 20000
$ws.Range("M54").Value = -19514
$ws.Range("H62").Value = 3395.35
$ws.Range("I62").Value = 3271
$ws.Range("K62").Value = 3271
$ws.Range("M62").Value = -2647
$ws.Range("H65").Value = 3395.35
$ws.Range("I65").Value = 3271
$ws.Range("K65").Value = 16355
$ws.Range("M65").Value = -13235
$ws.Range("H92").Value = 804.5909
$ws.Range("I92").Value = 680.05554
$ws.Range("J92").Value = 1365
$ws.Range("K92").Value = 680.05554
$ws.Range("L92").Value = 1365
$ws.Range("M92").Value = 567.94446
$ws.Range("N92").Value = -3861
$ws.Range("H103").Value = 400
$ws.Range("J103").Value = 400
$ws.Range("L103").Value = 1200
$ws.Range("N103").Value = -2372
$ws.Range("H130").Value = 19998.334
$ws.Range("J130").Value = 19998.334
$ws.Range("L130").Value = 19998.334
$ws.Range("N130").Value = -30038.334
$ws.Range("H132").Value = 12721.174
$ws.Range("I132").Value = 1513.6904
$ws.Range("K132").Value = 4541.0712
$ws.Range("M132").Value = -2011.0712
$ws.Range("H137").Value = 5188.1816
$ws.Range("I137").Value = 5430.1816
$ws.Range("K137").Value = 16290.5448
$ws.Range("M137").Value = -13740.5448

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4131.841
$ws.Range("I32").Value = 4714.123
$ws.Range("K32").Value = 4714.123
$ws.Range("M32").Value = -4427.123
$ws.Range("H45").Value = 3524.875
$ws.Range("H74").Value = 1448.0588
$ws.Range("I74").Value = 1497.6428
$ws.Range("K74").Value = 1497.6428
$ws.Range("M74").Value = -623.6428000000001
$ws.Range("H77").Value = 1448.0588
$ws.Range("I77").Value = 1497.6428
$ws.Range("K77").Value = 7488.214
$ws.Range("M77").Value = -3120.214
$ws.Range("H132").Value = 23257708
$ws.Range("I132").Value = 27779480
$ws.Range("K132").Value = 83338440
$ws.Range("M132").Value = -83335910

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1089.1666
$ws.Range("I36").Value = 1089.1666
$ws.Range("K36").Value = 1089.1666
$ws.Range("M36").Value = -555.1666
$ws.Range("H86").Value = 20834576
$ws.Range("I86").Value = 23810612
$ws.Range("J86").Value = 2324.3333
$ws.Range("K86").Value = 23810612
$ws.Range("L86").Value = 2324.3333
$ws.Range("M86").Value = -23809489
$ws.Range("N86").Value = -4570.3333
$ws.Range("H89").Value = 20834576
$ws.Range("I89").Value = 23810612
$ws.Range("J89").Value = 2324.3333
$ws.Range("K89").Value = 119053060
$ws.Range("L89").Value = 11621.6665
$ws.Range("M89").Value = -119047444
$ws.Range("N89").Value = -22853.6665
$ws.Range("H105").Value = 2726.5
$ws.Range("I105").Value = 1077
$ws.Range("J105").Value = 5200.75
$ws.Range("K105").Value = 1077
$ws.Range("L105").Value = 5200.75
$ws.Range("M105").Value = 670
$ws.Range("N105").Value = -8694.75
$ws.Range("H130").Value = 39428.57
$ws.Range("I130").Value = 0
$ws.Range("K130").Value = 0
$ws.Range("M130").ClearContents()
$ws.Range("H134").Value = 4750
$ws.Range("I134").Value = 4750
$ws.Range("K134").Value = 14250
$ws.Range("M134").Value = -11715

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1723.4
$ws.Range("I58").Value = 1609
$ws.Range("J58").Value = 1837.8
$ws.Range("K58").Value = 1609
$ws.Range("L58").Value = 1837.8
$ws.Range("M58").Value = -1406
$ws.Range("N58").Value = -2243.8
$ws.Range("H62").Value = 45458344
$ws.Range("J62").Value = 83337050
$ws.Range("L62").Value = 83337050
$ws.Range("N62").Value = -83338298
$ws.Range("H65").Value = 45458344
$ws.Range("J65").Value = 83337050
$ws.Range("L65").Value = 416685250
$ws.Range("N65").Value = -416691490
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H81").Value = 85000
$ws.Range("J81").Value = 85000
$ws.Range("L81").Value = 85000
$ws.Range("N81").Value = -86996
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H84").Value = 85000
$ws.Range("J84").Value = 85000
$ws.Range("L84").Value = 255000
$ws.Range("N84").Value = -264984
$ws.Range("H118").Value = 29833.334
$ws.Range("J118").Value = 29833.334
$ws.Range("L118").Value = 29833.334
$ws.Range("N118").Value = -33147.334
$ws.Range("H122").Value = 28064716
$ws.Range("I122").Value = 40410396
$ws.Range("K122").Value = 121231188
$ws.Range("M122").Value = -121228738
$ws.Range("H132").Value = 3237.15
$ws.Range("I132").Value = 2234.25
$ws.Range("K132").Value = 6702.75
$ws.Range("M132").Value = -4172.75
$ws.Range("H136").Value = 1723.4
$ws.Range("I136").Value = 1609
$ws.Range("J136").Value = 1837.8
$ws.Range("K136").Value = 4827
$ws.Range("L136").Value = 5513.4
$ws.Range("M136").Value = -2277
$ws.Range("N136").Value = -10613.4

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6655.7334
$ws.Range("I3").Value = 5702.5713
$ws.Range("K3").Value = 17107.7139
$ws.Range("M3").Value = -16995.7139
$ws.Range("H4").Value = 67000396
$ws.Range("I4").Value = 67000396
$ws.Range("K4").Value = 201001188
$ws.Range("M4").Value = -201001076
$ws.Range("H17").Value = 62.5
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H95").Value = 4930
$ws.Range("I95").Value = 4890
$ws.Range("J95").Value = 4950
$ws.Range("K95").Value = 14670
$ws.Range("L95").Value = 14850
$ws.Range("M95").Value = -12611
$ws.Range("N95").Value = -18968

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 3300
$ws.Range("J18").Value = 5000
$ws.Range("L18").Value = 5000
$ws.Range("N18").Value = -5586
$ws.Range("H80").Value = 99092.836
$ws.Range("I80").Value = 282277.25
$ws.Range("J80").Value = 7500.625
$ws.Range("K80").Value = 282277.25
$ws.Range("L80").Value = 7500.625
$ws.Range("M80").Value = -281279.25
$ws.Range("N80").Value = -9496.625
$ws.Range("H83").Value = 99092.836
$ws.Range("I83").Value = 282277.25
$ws.Range("J83").Value = 7500.625
$ws.Range("K83").Value = 1411386.25
$ws.Range("L83").Value = 37503.125
$ws.Range("M83").Value = -1406394.25
$ws.Range("N83").Value = -47487.125

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2707.4666
$ws.Range("I46").Value = 1996.5
$ws.Range("K46").Value = 1996.5
$ws.Range("M46").Value = -1808.5
$ws.Range("H68").Value = 1066.6666
$ws.Range("J68").Value = 1200
$ws.Range("L68").Value = 1200
$ws.Range("N68").Value = -2698
$ws.Range("H71").Value = 1066.6666
$ws.Range("J71").Value = 1200
$ws.Range("L71").Value = 6000
$ws.Range("N71").Value = -13488
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 3599.5
$ws.Range("I8").Value = 2199.5
$ws.Range("J8").Value = 4999.5
$ws.Range("K8").Value = 2199.5
$ws.Range("L8").Value = 4999.5
$ws.Range("M8").Value = -2059.5
$ws.Range("N8").Value = -5279.5
$ws.Range("H132").Value = 2817.7646
$ws.Range("I132").Value = 2673.4783
$ws.Range("K132").Value = 8020.4349
$ws.Range("M132").Value = -5490.4349
$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280
